$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.675.98"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "1.887.38"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "235.55"
$ws.Range("E5").Value = "  -4.06%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "0.4881"
$ws.Range("E7").Value = "  -2.09%  "
$ws.Range("D8").Value = "0.2880"
$ws.Range("E8").Value = "  -4.03%  "
$ws.Range("D9").Value = "0.06662"
$ws.Range("E9").Value = "  -3.45%  "
$ws.Range("D10").Value = "1.879.38"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").Value = "16.79"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").Value = "0.07240"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "88.75"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "5.006"
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("D15").Value = "0.6637"
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("D16").Value = "30.602.04"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "0.000007840"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "12.95"
$ws.Range("E19").Value = "  -3.65%  "
$ws.Range("D20").Value = "2.120.50"
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "4.734"
$ws.Range("E22").Value = "  -2.99%  "
$ws.Range("D23").Value = "186.68"
$ws.Range("E23").Value = "  +6.19%  "
$ws.Range("D24").Value = "6.040"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "157.84"
$ws.Range("E26").Value = "  +4.05%  "
$ws.Range("D27").Value = "18.28"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").Value = "1.831"
$ws.Range("E28").Value = "  -6.06%  "
$ws.Range("D29").Value = "1.407"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "4.256"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").Value = "0.09020"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").Value = "3.933"
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("D33").Value = "0.05196"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").Value = "0.7314"
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("D35").Value = "1.079"
$ws.Range("E35").Value = "  -5.71%  "
$ws.Range("D36").Value = "2.700"
$ws.Range("D37").Value = "0.01817"
$ws.Range("E37").Value = "  -5.00%  "
$ws.Range("D38").Value = "2.657"
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").Value = "0.9189"
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("E40").Value = "  -7.49%  "
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("D42").Value = "104.05"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").Value = "0.9985"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "5.716"
$ws.Range("E44").Value = "  -3.85%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").Value = "7.275"
$ws.Range("E46").Value = "  -6.93%  "
$ws.Range("D47").Value = "0.05831"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.3945"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.619"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").Value = "1.410"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("D51").Value = "33.20"
$ws.Range("E51").Value = "  -0.48%  "